# TST: Test unnamed columns with index_col for Excel (#23874)
# Adds a new Sheet4 to the workbook that exercises a DataFrame written with
# an index column and unnamed header columns - mirroring pandas'
# ExcelFormatter default header style (bold, thin-bordered, centered).

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet -------------------------------------------------
$ws = $wb.Worksheets.Add()
$ws.Name = "Sheet4"

# --- Cell values ---------------------------------------------------------------
# Row 1: unnamed-index header row -> columns "col1"/"col2" start at column B
$ws.Range("B1").Value = "col1"
$ws.Range("C1").Value = "col2"

# Row 2: first data row, index label "i1" in column A
$ws.Range("A2").Value = "i1"
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "x"

# Row 3: second data row, index label "i2" in column A
$ws.Range("A3").Value = "i2"
$ws.Range("B3").Value = "b"
$ws.Range("C3").Value = "y"

# --- Header styling (bold Calibri 11, thin border box, center/top align) -------
function Set-HeaderStyle($rng) {
    $f = $rng.Font
    $f.Name = "Calibri"
    $f.Size = 11
    $f.Bold = $true
    $f.ThemeColor = 1
    $f.ThemeFont = 1

    $rng.Borders.LineStyle = 1

    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

Set-HeaderStyle $ws.Range("B1:C1")
Set-HeaderStyle $ws.Range("A2:A3")

# --- Row heights to match the header font metrics -------------------------------
$ws.Rows.Item(1).RowHeight = 14.4
$ws.Rows.Item(2).RowHeight = 14.4
$ws.Rows.Item(3).RowHeight = 14.4

# --- Move the new sheet to the end of the tab strip and select its full range --
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A1:C3").Select()
$ws4.Activate()
